# Add a new "Italy" test-data sheet, cloned from "Slovakia", with its own
# Market / NGC values, inserted right after "Slovakia" and made the active
# (selected) sheet/cell — matching the author's "Test data added for Italy"
# commit.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy Slovakia and place the copy immediately after the last sheet
# (i.e. right after Slovakia, since Slovakia currently is the last sheet).
$slovakia.Copy($null, $lastSheet)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Fill in the country-specific values for the new sheet.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2224 "

# Restore Slovakia's view to a plain "select all" state (no longer the
# active tab) and make Italy the active tab with B4 selected.
$slovakia.Range("A1:XFD1048576").Select()

$italy.Select()
$italy.Range("B4").Select()
